$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 757.61536
$ws.Cells.Item(41, 9).Value = 641
$ws.Cells.Item(41, 10).Value = 944.2
$ws.Cells.Item(41, 11).Value = 641
$ws.Cells.Item(41, 12).Value = 944.2
$ws.Cells.Item(41, 13).Value = -201
$ws.Cells.Item(41, 14).Value = -1824.2
$ws.Cells.Item(43, 8).Value = 7872.7144
$ws.Cells.Item(43, 9).Value = 12562.5
$ws.Cells.Item(43, 10).Value = 5996.8
$ws.Cells.Item(43, 11).Value = 12562.5
$ws.Cells.Item(43, 12).Value = 5996.8
$ws.Cells.Item(43, 13).Value = -12493.5
$ws.Cells.Item(43, 14).Value = -6134.8
$ws.Cells.Item(64, 8).Value = 6493.5713
$ws.Cells.Item(64, 9).Value = 3455
$ws.Cells.Item(64, 10).Value = 7000
$ws.Cells.Item(64, 11).Value = 3455
$ws.Cells.Item(64, 12).Value = 7000
$ws.Cells.Item(64, 13).Value = -3207
$ws.Cells.Item(64, 14).Value = -7496
$ws.Cells.Item(67, 8).Value = 6493.5713
$ws.Cells.Item(67, 9).Value = 3455
$ws.Cells.Item(67, 10).Value = 7000
$ws.Cells.Item(67, 11).Value = 3455
$ws.Cells.Item(67, 12).Value = 7000
$ws.Cells.Item(67, 13).Value = -2597
$ws.Cells.Item(67, 14).Value = -8716
$ws.Cells.Item(108, 8).Value = 81995
$ws.Cells.Item(108, 10).Value = 81995
$ws.Cells.Item(108, 12).Value = 81995
$ws.Cells.Item(108, 14).Value = -89675
$ws.Cells.Item(137, 8).Value = 26602.176
$ws.Cells.Item(137, 9).Value = 30812.896
$ws.Cells.Item(137, 10).Value = 2180
$ws.Cells.Item(137, 11).Value = 92438.68799999999
$ws.Cells.Item(137, 12).Value = 6540
$ws.Cells.Item(137, 13).Value = -89888.68799999999
$ws.Cells.Item(137, 14).Value = -11640
$ws.Cells.Item(138, 8).Value = 5295.5264
$ws.Cells.Item(138, 9).Value = 1806.8182
$ws.Cells.Item(138, 11).Value = 5420.4546
$ws.Cells.Item(138, 13).Value = -280.4546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4254229.5
$ws.Cells.Item(32, 9).Value = 1794462.2
$ws.Cells.Item(32, 11).Value = 1794462.2
$ws.Cells.Item(32, 13).Value = -1794175.2
$ws.Cells.Item(45, 8).Value = 5841.8237
$ws.Cells.Item(45, 9).Value = 4351.5
$ws.Cells.Item(45, 10).Value = 7166.5557
$ws.Cells.Item(45, 11).Value = 4351.5
$ws.Cells.Item(45, 12).Value = 7166.5557
$ws.Cells.Item(45, 13).Value = -3974.5
$ws.Cells.Item(45, 14).Value = -7920.5557
$ws.Cells.Item(61, 8).Value = 3934.5293
$ws.Cells.Item(61, 9).Value = 3908.0908
$ws.Cells.Item(61, 11).Value = 3908.0908
$ws.Cells.Item(61, 13).Value = -3696.0908
$ws.Cells.Item(122, 8).Value = 5059.2334
$ws.Cells.Item(122, 9).Value = 3698.5715
$ws.Cells.Item(122, 10).Value = 6249.8125
$ws.Cells.Item(122, 11).Value = 11095.7145
$ws.Cells.Item(122, 12).Value = 18749.4375
$ws.Cells.Item(122, 13).Value = -8645.7145
$ws.Cells.Item(122, 14).Value = -23649.4375
$ws.Cells.Item(136, 8).Value = 3934.5293
$ws.Cells.Item(136, 9).Value = 3908.0908
$ws.Cells.Item(136, 11).Value = 11724.2724
$ws.Cells.Item(136, 13).Value = -9174.2724

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1200.8966
$ws.Cells.Item(107, 9).Value = 1166.2084
$ws.Cells.Item(107, 11).Value = 1166.2084
$ws.Cells.Item(107, 13).Value = 753.7916
$ws.Cells.Item(134, 8).Value = 8931221
$ws.Cells.Item(134, 9).Value = 14288154
$ws.Cells.Item(134, 10).Value = 3000
$ws.Cells.Item(134, 11).Value = 42864462
$ws.Cells.Item(134, 12).Value = 9000
$ws.Cells.Item(134, 13).Value = -42861927
$ws.Cells.Item(134, 14).Value = -14070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2805.4075
$ws.Cells.Item(58, 9).Value = 2715.861
$ws.Cells.Item(58, 11).Value = 2715.861
$ws.Cells.Item(58, 13).Value = -2512.861
$ws.Cells.Item(62, 8).Value = 3200.6
$ws.Cells.Item(62, 10).Value = 3333.3333
$ws.Cells.Item(62, 12).Value = 3333.3333
$ws.Cells.Item(62, 14).Value = -4581.3333
$ws.Cells.Item(65, 8).Value = 3200.6
$ws.Cells.Item(65, 10).Value = 3333.3333
$ws.Cells.Item(65, 12).Value = 16666.6665
$ws.Cells.Item(65, 14).Value = -22906.6665
$ws.Cells.Item(122, 8).Value = 1621.125
$ws.Cells.Item(122, 9).Value = 1694.2
$ws.Cells.Item(122, 11).Value = 5082.6
$ws.Cells.Item(122, 13).Value = -2632.6
$ws.Cells.Item(136, 8).Value = 2805.4075
$ws.Cells.Item(136, 9).Value = 2715.861
$ws.Cells.Item(136, 11).Value = 8147.583
$ws.Cells.Item(136, 13).Value = -5597.583
$ws.Cells.Item(141, 8).Value = 360239.62
$ws.Cells.Item(141, 10).Value = 360239.62
$ws.Cells.Item(141, 12).Value = 360239.62
$ws.Cells.Item(141, 14).Value = -370599.62

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 101157530
$ws.Cells.Item(4, 10).Value = 66665720
$ws.Cells.Item(4, 12).Value = 199997160
$ws.Cells.Item(4, 14).Value = -199997384
$ws.Cells.Item(11, 8).Value = 573311.1
$ws.Cells.Item(11, 9).Value = 1010760.4
$ws.Cells.Item(11, 10).Value = 26499.5
$ws.Cells.Item(11, 11).Value = 3032281.2
$ws.Cells.Item(11, 12).Value = 79498.5
$ws.Cells.Item(11, 13).Value = -3032141.2
$ws.Cells.Item(11, 14).Value = -79778.5
$ws.Cells.Item(12, 8).Value = 14.571428
$ws.Cells.Item(12, 9).Value = 15.666667
$ws.Cells.Item(12, 10).Value = 13.75
$ws.Cells.Item(12, 11).Value = 47.000001
$ws.Cells.Item(12, 12).Value = 41.25
$ws.Cells.Item(12, 13).Value = 125.999999
$ws.Cells.Item(12, 14).Value = -387.25
$ws.Cells.Item(13, 8).Value = 163.33333
$ws.Cells.Item(13, 10).Value = 195
$ws.Cells.Item(13, 12).Value = 585
$ws.Cells.Item(13, 14).Value = -921
$ws.Cells.Item(38, 8).Value = 417.25
$ws.Cells.Item(38, 9).Value = 62.5
$ws.Cells.Item(38, 10).Value = 653.75
$ws.Cells.Item(38, 11).Value = 187.5
$ws.Cells.Item(38, 12).Value = 1961.25
$ws.Cells.Item(38, 13).Value = 159.5
$ws.Cells.Item(38, 14).Value = -2655.25
$ws.Cells.Item(39, 8).Value = 3696.5
$ws.Cells.Item(39, 10).Value = 4927.6665
$ws.Cells.Item(39, 12).Value = 14782.9995
$ws.Cells.Item(39, 14).Value = -15370.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 289
$ws.Cells.Item(2, 10).Value = 556.3333
$ws.Cells.Item(2, 12).Value = 556.3333
$ws.Cells.Item(2, 14).Value = -782.3333
$ws.Cells.Item(5, 8).Value = 9996.143
$ws.Cells.Item(5, 9).Value = 9996.143
$ws.Cells.Item(5, 11).Value = 9996.143
$ws.Cells.Item(5, 13).Value = -9884.143
$ws.Cells.Item(122, 8).Value = 1419.8334
$ws.Cells.Item(122, 9).Value = 1419.8334
$ws.Cells.Item(122, 11).Value = 4259.5002
$ws.Cells.Item(122, 13).Value = -1809.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 377.1905
$ws.Cells.Item(55, 9).Value = 235.63637
$ws.Cells.Item(55, 10).Value = 532.9
$ws.Cells.Item(55, 11).Value = 235.63637
$ws.Cells.Item(55, 12).Value = 532.9
$ws.Cells.Item(55, 13).Value = -62.63637
$ws.Cells.Item(55, 14).Value = -878.9
$ws.Cells.Item(68, 8).Value = 3667.3333
$ws.Cells.Item(68, 10).Value = 3500
$ws.Cells.Item(68, 12).Value = 3500
$ws.Cells.Item(68, 14).Value = -4998
$ws.Cells.Item(71, 8).Value = 3667.3333
$ws.Cells.Item(71, 10).Value = 3500
$ws.Cells.Item(71, 12).Value = 17500
$ws.Cells.Item(71, 14).Value = -24988
$ws.Cells.Item(93, 8).Value = 71430280
$ws.Cells.Item(93, 9).Value = 111112550
$ws.Cells.Item(93, 10).Value = 2182.2
$ws.Cells.Item(93, 11).Value = 111112550
$ws.Cells.Item(93, 12).Value = 2182.2
$ws.Cells.Item(93, 13).Value = -111111302
$ws.Cells.Item(93, 14).Value = -4678.2
$ws.Cells.Item(122, 8).Value = 9973.223
$ws.Cells.Item(122, 9).Value = 8682.200000000001
$ws.Cells.Item(122, 11).Value = 26046.6
$ws.Cells.Item(122, 13).Value = -23596.6
$ws.Cells.Item(136, 8).Value = 4247.136
$ws.Cells.Item(136, 9).Value = 2948.4546
$ws.Cells.Item(136, 11).Value = 8845.363799999999
$ws.Cells.Item(136, 13).Value = -6295.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 13129708
$ws.Cells.Item(2, 9).Value = 14291095
$ws.Cells.Item(2, 10).Value = 5000000
$ws.Cells.Item(2, 11).Value = 14291095
$ws.Cells.Item(2, 12).Value = 5000000
$ws.Cells.Item(2, 13).Value = -14290983
$ws.Cells.Item(2, 14).Value = -5000224
$ws.Cells.Item(122, 8).Value = 4031.075
$ws.Cells.Item(122, 9).Value = 2269.9119
$ws.Cells.Item(122, 11).Value = 6809.7357
$ws.Cells.Item(122, 13).Value = -4359.7357
